$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every student row (2..50), if the "total_views" (column I) is non-zero,
# zero out the daily view flags (columns B..H) and the computed totals
# (I: total_views, J: nota_view). Rows already all-zero are left untouched.
for ($row = 2; $row -le 50; $row++) {
    $totalViews = $ws.Cells.Item($row, 9).Value
    if ($totalViews -ne 0) {
        for ($col = 2; $col -le 10; $col++) {
            $ws.Cells.Item($row, $col).Value = 0
        }
    }
}
